$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range('D2:E2')
$rng.NumberFormat = "@"
$ws.Range('D2').Value = '315.81'
$ws.Range('E2').Value = '6.54%'
$rng.Style = "Normal"

$rng = $ws.Range('D3:E3')
$rng.NumberFormat = "@"
$ws.Range('D3').Value = '45.28'
$ws.Range('E3').Value = '7.47%'
$rng.Style = "Normal"

$rng = $ws.Range('D4:E4')
$rng.NumberFormat = "@"
$ws.Range('D4').Value = '5.177'
$ws.Range('E4').Value = '2.93%'
$rng.Style = "Normal"

$rng = $ws.Range('D5:E5')
$rng.NumberFormat = "@"
$ws.Range('D5').Value = '0.08083'
$ws.Range('E5').Value = '6.54%'
$rng.Style = "Normal"

$rng = $ws.Range('E6')
$rng.NumberFormat = "@"
$ws.Range('E6').Value = '3.21%'
$rng.Style = "Normal"

$rng = $ws.Range('D7:E7')
$rng.NumberFormat = "@"
$ws.Range('D7').Value = '1.670'
$ws.Range('E7').Value = '4.05%'
$rng.Style = "Normal"

$rng = $ws.Range('D8:E8')
$rng.NumberFormat = "@"
$ws.Range('D8').Value = '1.093'
$ws.Range('E8').Value = '17.43%'
$rng.Style = "Normal"

$rng = $ws.Range('D9:E9')
$rng.NumberFormat = "@"
$ws.Range('D9').Value = '0.1307'
$ws.Range('E9').Value = '8.20%'
$rng.Style = "Normal"

$rng = $ws.Range('D10:E10')
$rng.NumberFormat = "@"
$ws.Range('D10').Value = '0.1935'
$ws.Range('E10').Value = '5.09%'
$rng.Style = "Normal"

$rng = $ws.Range('D11:E11')
$rng.NumberFormat = "@"
$ws.Range('D11').Value = '0.09504'
$ws.Range('E11').Value = '5.51%'
$rng.Style = "Normal"

$rng = $ws.Range('D12:E12')
$rng.NumberFormat = "@"
$ws.Range('D12').Value = '0.04234'
$ws.Range('E12').Value = '6.14%'
$rng.Style = "Normal"

$rng = $ws.Range('E13')
$rng.NumberFormat = "@"
$ws.Range('E13').Value = '-0.83%'
$rng.Style = "Normal"

$rng = $ws.Range('D14:E14')
$rng.NumberFormat = "@"
$ws.Range('D14').Value = '0.001315'
$ws.Range('E14').Value = '2.79%'
$rng.Style = "Normal"

$rng = $ws.Range('D15:E15')
$rng.NumberFormat = "@"
$ws.Range('D15').Value = '0.005938'
$ws.Range('E15').Value = '1.28%'
$rng.Style = "Normal"

$rng = $ws.Range('E17')
$rng.NumberFormat = "@"
$ws.Range('E17').Value = '1.03%'
$rng.Style = "Normal"

$rng = $ws.Range('D18:E18')
$rng.NumberFormat = "@"
$ws.Range('D18').Value = '2.410'
$ws.Range('E18').Value = '0.20%'
$rng.Style = "Normal"

$rng = $ws.Range('D19:E19')
$rng.NumberFormat = "@"
$ws.Range('D19').Value = '0.3370'
$ws.Range('E19').Value = '1.50%'
$rng.Style = "Normal"

$rng = $ws.Range('D20:E20')
$rng.NumberFormat = "@"
$ws.Range('D20').Value = '8.207'
$ws.Range('E20').Value = '4.10%'
$rng.Style = "Normal"

$rng = $ws.Range('D21:E21')
$rng.NumberFormat = "@"
$ws.Range('D21').Value = '0.1385'
$ws.Range('E21').Value = '-2.41%'
$rng.Style = "Normal"

$rng = $ws.Range('D23:E23')
$rng.NumberFormat = "@"
$ws.Range('D23').Value = '0.04283'
$ws.Range('E23').Value = '5.45%'
$rng.Style = "Normal"

$rng = $ws.Range('E24')
$rng.NumberFormat = "@"
$ws.Range('E24').Value = '1.20%'
$rng.Style = "Normal"

$rng = $ws.Range('D25:E25')
$rng.NumberFormat = "@"
$ws.Range('D25').Value = '0.004225'
$ws.Range('E25').Value = '7.85%'
$rng.Style = "Normal"

$rng = $ws.Range('D26:E26')
$rng.NumberFormat = "@"
$ws.Range('D26').Value = '0.0001346'
$ws.Range('E26').Value = '9.40%'
$rng.Style = "Normal"

$rng = $ws.Range('D38:E38')
$rng.NumberFormat = "@"
$ws.Range('D38').Value = '0.02711'
$ws.Range('E38').Value = '11.93%'
$rng.Style = "Normal"

$rng = $ws.Range('E39')
$rng.NumberFormat = "@"
$ws.Range('E39').Value = '4.81%'
$rng.Style = "Normal"

$rng = $ws.Range('D40:E40')
$rng.NumberFormat = "@"
$ws.Range('D40').Value = '0.005868'
$ws.Range('E40').Value = '-3.05%'
$rng.Style = "Normal"

$rng = $ws.Range('D41:E41')
$rng.NumberFormat = "@"
$ws.Range('D41').Value = '0.007778'
$ws.Range('E41').Value = '-0.15%'
$rng.Style = "Normal"

$rng = $ws.Range('D42:E42')
$rng.NumberFormat = "@"
$ws.Range('D42').Value = '0.1426'
$ws.Range('E42').Value = '7.03%'
$rng.Style = "Normal"

$rng = $ws.Range('D43:E43')
$rng.NumberFormat = "@"
$ws.Range('D43').Value = '0.007374'
$ws.Range('E43').Value = '-2.25%'
$rng.Style = "Normal"

$rng = $ws.Range('D44:E44')
$rng.NumberFormat = "@"
$ws.Range('D44').Value = '0.008597'
$ws.Range('E44').Value = '18.55%'
$rng.Style = "Normal"

$rng = $ws.Range('D45:E45')
$rng.NumberFormat = "@"
$ws.Range('D45').Value = '0.3144'
$ws.Range('E45').Value = '5.67%'
$rng.Style = "Normal"

$rng = $ws.Range('D46:E46')
$rng.NumberFormat = "@"
$ws.Range('D46').Value = '0.00006803'
$ws.Range('E46').Value = '0.20%'
$rng.Style = "Normal"

$rng = $ws.Range('D48:E48')
$rng.NumberFormat = "@"
$ws.Range('D48').Value = '0.06227'
$ws.Range('E48').Value = '35.78%'
$rng.Style = "Normal"

$rng = $ws.Range('D49:E49')
$rng.NumberFormat = "@"
$ws.Range('D49').Value = '0.003985'
$ws.Range('E49').Value = '-5.17%'
$rng.Style = "Normal"

$rng = $ws.Range('D51')
$rng.NumberFormat = "@"
$ws.Range('D51').Value = '0.0001995'
$rng.Style = "Normal"
